$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data / mean calculation
$ws.Range("F3").Value = -4
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = -6
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = -5
